# Changes to improve performance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C2:C11 to use the new checkbox-4 identifier (a new shared string
# is introduced; the old checkbox-3 string at C1 area is left untouched).
$ws.Range("C2:C11").Value = "individualQuotaValuesId-row-checkbox-4"

# Update the active selection on Sheet1 to C3:C11 with active cell C3.
$ws.Range("C3:C11").Select()
